$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 227-228, shifting the existing data (previously
# rows 227-284) down to rows 229-286. This matches the diff, where every
# existing record moved down by two rows and two brand-new weekly price
# records were introduced at the top of this block.
$ws.Range("A227:A228").EntireRow.Insert()

# New row 227: Apio, Americana (o), Primera - week of 2022-06-24
$ws.Cells.Item(227, 1).Value = 11
$ws.Cells.Item(227, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(227, 3).Value = "Bíobío"
$ws.Cells.Item(227, 4).Value = 44736
$ws.Cells.Item(227, 5).Value = 8
$ws.Cells.Item(227, 6).Value = 100112017
$ws.Cells.Item(227, 7).Value = "Apio"
$ws.Cells.Item(227, 8).Value = "Americana (o)"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 100
$ws.Cells.Item(227, 11).Value = 7000
$ws.Cells.Item(227, 12).Value = 7500
$ws.Cells.Item(227, 13).Value = 7250
$ws.Cells.Item(227, 14).Value = "`$/docena de matas"
$ws.Cells.Item(227, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(227, 16).Value = 1208
$ws.Cells.Item(227, 17).Value = 6
$ws.Cells.Item(227, 18).Value = "Hortaliza"

# New row 228: Apio, Americana (o), Segunda - week of 2022-06-24
$ws.Cells.Item(228, 1).Value = 11
$ws.Cells.Item(228, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(228, 3).Value = "Bíobío"
$ws.Cells.Item(228, 4).Value = 44736
$ws.Cells.Item(228, 5).Value = 8
$ws.Cells.Item(228, 6).Value = 100112017
$ws.Cells.Item(228, 7).Value = "Apio"
$ws.Cells.Item(228, 8).Value = "Americana (o)"
$ws.Cells.Item(228, 9).Value = "Segunda"
$ws.Cells.Item(228, 10).Value = 50
$ws.Cells.Item(228, 11).Value = 6500
$ws.Cells.Item(228, 12).Value = 6500
$ws.Cells.Item(228, 13).Value = 6500
$ws.Cells.Item(228, 14).Value = "`$/docena de matas"
$ws.Cells.Item(228, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(228, 16).Value = 1083
$ws.Cells.Item(228, 17).Value = 6
$ws.Cells.Item(228, 18).Value = "Hortaliza"
